$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text format to avoid numeric auto-conversion
$priceCells = @("D2", "D3", "D5", "D8", "D11", "D12", "D13", "D14", "D15", "D18", "D19", "D22", "D25", "D26", "D27", "D30", "D33", "D35", "D36", "D40", "D43", "D44", "D47", "D48", "D50")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.128.50"
$ws.Range("D3").Value = "1.781.87"
$ws.Range("D5").Value = "226.02"
$ws.Range("D8").Value = "32.13"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D12").Value = "2.038.75"
$ws.Range("D13").Value = "1.784.22"
$ws.Range("D14").Value = "10.95"
$ws.Range("D15").Value = "34.115.28"
$ws.Range("D18").Value = "67.58"
$ws.Range("D19").Value = "245.39"
$ws.Range("D22").Value = "10.88"
$ws.Range("D25").Value = "161.93"
$ws.Range("D26").Value = "7.13"
$ws.Range("D27").Value = "16.28"
$ws.Range("D30").Value = "1.22"
$ws.Range("D33").Value = "3.72"
$ws.Range("D35").Value = "1.443.68"
$ws.Range("D36").Value = "2.49"
$ws.Range("D40").Value = "81.43"
$ws.Range("D43").Value = "0.913"
$ws.Range("D44").Value = "13.61"
$ws.Range("D47").Value = "6.05"
$ws.Range("D48").Value = "1.938.76"
$ws.Range("D50").Value = "104.75"

foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("E14").Value = "  -4.35%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("E36").Value = "  +6.43%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  -6.53%  "
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("E51").Value = "  +0.25%  "
